# "Generate Report for Handback"
#
# Populates the "Latest Target File" (F) and "Latest Handback File" (G)
# columns -- previously empty -- on the per-language handoff sheets, flips
# the status from "Ready for handoff" to "Handed back: in sync with en-US",
# and stamps the "Latest Handback DateTime" (H) column with the actual
# handback timestamps (replacing the zero-date placeholder).

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdName  = "0e71f548-af60-45c4-b839-55682e57a859.md"
$zhXlf   = "0e71f548-af60-45c4-b839-55682e57a859.0aa05150851e84d6a81556a899c19be5e8830348.zh-cn.xlf"
$deXlf   = "0e71f548-af60-45c4-b839-55682e57a859.0aa05150851e84d6a81556a899c19be5e8830348.de-de.xlf"

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/bc5727dcb531dcfb241abb13bb74d469146173b0/e2e/" + $mdName
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b9393c899e49546261c1c0135d441e0309bbb078/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $zhXlf
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac4ccad78c9b2ee56599052e0680ab03f32aedc7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $deXlf

# --- Overview sheet: status text only (columns B/C, rows 2-3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusNew
$zh.Range("C3").Value = $statusNew

$zh.Hyperlinks.Add($zh.Range("F2"), $mdUrl, "", "", $mdName)
$zh.Hyperlinks.Add($zh.Range("G2"), $zhXlfUrl, "", "", $zhXlf)
$zh.Hyperlinks.Add($zh.Range("F3"), $mdUrl, "", "", $mdName)
$zh.Hyperlinks.Add($zh.Range("G3"), $zhXlfUrl, "", "", $zhXlf)

$zh.Range("H2").Value = "2016-03-11 09:47:01"
$zh.Range("H3").Value = "2016-03-11 09:47:01"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusNew
$de.Range("C3").Value = $statusNew

$de.Hyperlinks.Add($de.Range("F2"), $mdUrl, "", "", $mdName)
$de.Hyperlinks.Add($de.Range("G2"), $deXlfUrl, "", "", $deXlf)
$de.Hyperlinks.Add($de.Range("F3"), $mdUrl, "", "", $mdName)
$de.Hyperlinks.Add($de.Range("G3"), $deXlfUrl, "", "", $deXlf)

$de.Range("H2").Value = "2016-03-11 09:47:07"
$de.Range("H3").Value = "2016-03-11 09:47:07"

Write-Host "Handback report generated."
